# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets, matching the freshly generated data
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F on sheet "展览"
$exhibitionUpdates = @{
    6  = 2837
    8  = 1846
    10 = 89
    11 = 667
    14 = 195
    16 = 90
    17 = 27
}

# Row -> new value for column F on sheet "全部类型"
$allTypesUpdates = @{
    6  = 2837
    8  = 1846
    10 = 89
    11 = 667
    14 = 195
    16 = 90
    17 = 27
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
